$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "StatQuery" Cypher query text (shared by C2, C3, C4) ---
$newQuery = "MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)`nOPTIONAL MATCH (samp:sample)-->(c)`nOPTIONAL MATCH (diag:diagnosis)-->(c)`nOPTIONAL MATCH (f:file)-[*]->(c)`nOPTIONAL MATCH (sf:file)-->(s)`nWITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p`nWHERE demo.breed IN ['Australian Shepherd']`nRETURN  `n    count(distinct p) AS Programs,`n    count(distinct s) AS Studies,`n    count(distinct c) AS Cases,`n    count(distinct samp) AS Samples,`n    count(distinct f) AS ``Case Files``,`n    count(distinct sf) AS ``Study Files``"

$ws.Range("C2").Value = $newQuery
$ws.Range("C3").Value = $newQuery
$ws.Range("C4").Value = $newQuery

# --- Fix up the saved view state: scroll back to the top and select B1 ---
$ws.Range("B1").Select()
